$wb = $excel.ActiveWorkbook

# --- Sheet "Classes" ---
$ws1 = $wb.Worksheets.Item("Classes")

# New column I header (Utils class)
$ws1.Range("I1").Value = "Utils"

# Row 2 (Fields) updates
$ws1.Range("C2").Value = "Square[][] squares`nint roundNumber`nGamePhase phase`nPlayer winner"
$ws1.Range("E2").Value = "`nPlayer owner"
$ws1.Range("H2").Value = "Node treeRoot"

# Row 3 (Methods) updates
$ws1.Range("H3").Value = "void train(Node root, int seconds)`nvoid select(Node node, int Ni)`nvoid simulate(Node node)`nvoid backPropogate(Node node, Player winner)"
$ws1.Range("I3").Value = "Board getNextBoard(Board, Delta)"

# Row 5 - clear old note, keep its wrapped style
$ws1.Range("H5").ClearContents()

# Column width adjustments (G widened, new H column added)
$ws1.Columns.Item(7).ColumnWidth = 9.25
$ws1.Columns.Item(8).ColumnWidth = 44.3

# Row height for row 3 grew because of the new 4-line method list
$ws1.Rows.Item(3).RowHeight = 61.5

$ws1.Activate() | Out-Null
$ws1.Range("I2").Select() | Out-Null

# --- Sheet "Misc" ---
$ws3 = $wb.Worksheets.Item("Misc")
$ws3.Range("D1").Value = "Player"
$ws3.Range("D2").Value = "Enum"
$ws3.Range("D3").Value = "WHITE`nBLACK`n.getRepresentation() (@ or O)"
$ws3.Range("D3").WrapText = $true

$ws3.Activate() | Out-Null
$ws3.Range("D3").Select() | Out-Null

$ws1.Activate() | Out-Null
